$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that were removed entirely (naive forecaster bug fix removed some stale forecast rows)
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("C4").ClearContents()

# Update recomputed forecast values with corrected (bugfixed) figures
$ws.Range("E3").Value = 2.829537440099972
$ws.Range("E4").Value = 2.957479223435744
$ws.Range("C5").Value = 0.5917823527752386
$ws.Range("E5").Value = -1.985049937500005
$ws.Range("C6").Value = -0.01587181126745385
$ws.Range("E6").Value = 3.8351443707757
$ws.Range("C8").Value = -0.02256889165886955
$ws.Range("E8").Value = -0.6757980944263275
$ws.Range("C9").Value = -0.7518797681959066
$ws.Range("C10").Value = 0.09611428386595566
$ws.Range("C11").Value = -0.5765930039053124
$ws.Range("E12").Value = -0.971238541762387
$ws.Range("C13").Value = -0.07642926654479743
$ws.Range("E13").Value = -0.3994003999000184
$ws.Range("C14").Value = -0.001350220946472191
$ws.Range("E14").Value = 0.6008487920565075
$ws.Range("C15").Value = 0.9274109147535459
$ws.Range("E15").Value = 2.82953744009995
$ws.Range("E16").Value = -1.58998093318411
$ws.Range("E17").Value = -1.194610791899997
$ws.Range("C18").Value = -0.5761528471665334
$ws.Range("E18").Value = 0.4501721032283301
$ws.Range("E20").Value = -0.150175137493469
$ws.Range("C21").Value = -0.07666472728170559
$ws.Range("C23").Value = 0.1986438914956423
$ws.Range("E24").Value = -0.3994003999000073
$ws.Range("E26").Value = 2.372078088364704
$ws.Range("C27").Value = -0.7530239469328737
$ws.Range("E27").Value = -2.378486270399993
$ws.Range("C28").Value = -0.7283174404323023
$ws.Range("E28").Value = -2.378486270399993
$ws.Range("C29").Value = -0.2044553505917812
$ws.Range("C30").Value = -0.4278219446121501
$ws.Range("E30").Value = -2.378564786744752
$ws.Range("E31").Value = 2.421686529599998
$ws.Range("C32").Value = -0.03096525636256953
$ws.Range("E32").Value = 1.205410808099949
$ws.Range("C33").Value = -0.2555583584977206
$ws.Range("C34").Value = -1.026566979837429
$ws.Range("C35").Value = 1.98829857406233
$ws.Range("E35").Value = 5.718701441600027
$ws.Range("C36").Value = 1.135350354669407
$ws.Range("E36").Value = 5.718701441600027
$ws.Range("E39").Value = 0.4006004000999486
$ws.Range("C40").Value = 2.185449115957461
$ws.Range("C42").Value = 0.4636049209196802
$ws.Range("C45").Value = 0.5210077780289701
$ws.Range("C46").Value = 0.6216390921348403
$ws.Range("E46").Value = -1.097580983230539
$ws.Range("C47").Value = -0.5760755550525465
$ws.Range("C49").Value = -0.8507045154764525
$ws.Range("C50").Value = -0.6768900623516871
$ws.Range("C51").Value = 2.67797050805143
$ws.Range("E51").Value = 3.238605209600021
$ws.Range("E52").Value = -3.551690943900021
